$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (shared string index 224, cell A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 18:33"

# Country-name corrections: these rows had the wrong country label attached
# to their statistics; swap the labels so the figures line up with the right country.
$ws.Range("A60").Value = "Azerbaiyan"
$ws.Range("A61").Value = "Honduras"
$ws.Range("A65").Value = "Argelia"
$ws.Range("A66").Value = "Camerun"
$ws.Range("A107").Value = "Sri Lanka"
$ws.Range("A108").Value = "Mali"
$ws.Range("A128").Value = "Jordania"
$ws.Range("A129").Value = "Niger"

# Refresh the COVID-19 statistics (columns B-H) with the latest figures.
$ws.Range("B4").Value = 2400257
$ws.Range("C4").Value = 12104
$ws.Range("D4").Value = 1003712
$ws.Range("E4").Value = 1273697
$ws.Range("G4").Value = 238
$ws.Range("H4").Value = 122848

$ws.Range("B8").Value = 306210
$ws.Range("C8").Value = 921
$ws.Range("G8").Value = 280
$ws.Range("H8").Value = 42927

$ws.Range("B12").Value = 238833
$ws.Range("C12").Value = 113
$ws.Range("D12").Value = 184585
$ws.Range("E12").Value = 19573
$ws.Range("G12").Value = 18
$ws.Range("H12").Value = 34675

$ws.Range("B21").Value = 101902
$ws.Range("C21").Value = 265
$ws.Range("D21").Value = 64508
$ws.Range("E21").Value = 28941
$ws.Range("G21").Value = 17
$ws.Range("H21").Value = 8453

$ws.Range("D35").Value = 35995
$ws.Range("E35").Value = 6411

$ws.Range("B39").Value = 34502
$ws.Range("C39").Value = 1826
$ws.Range("D39").Value = 15753
$ws.Range("E39").Value = 17498
$ws.Range("G39").Value = 84
$ws.Range("H39").Value = 1251

$ws.Range("D54").Value = 11514
$ws.Range("E54").Value = 6583
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 134

$ws.Range("B60").Value = 13715
$ws.Range("C60").Value = 508
$ws.Range("D60").Value = 7503
$ws.Range("E60").Value = 6045
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 167

$ws.Range("B61").Value = 13356
$ws.Range("C61").Value = 584
$ws.Range("D61").Value = 1362
$ws.Range("E61").Value = 11599
$ws.Range("G61").Value = 32
$ws.Range("H61").Value = 395

$ws.Range("B65").Value = 12077
$ws.Range("C65").Value = 157
$ws.Range("D65").Value = 8653
$ws.Range("E65").Value = 2563
$ws.Range("G65").Value = 9
$ws.Range("H65").Value = 861

$ws.Range("B66").Value = 12041
$ws.Range("D66").Value = 7740
$ws.Range("E66").Value = 3993
$ws.Range("H66").Value = 308

$ws.Range("B67").Value = 10619
$ws.Range("C67").Value = 96
$ws.Range("D67").Value = 7555
$ws.Range("E67").Value = 2725
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 339

$ws.Range("B88").Value = 4133
$ws.Range("C88").Value = 12
$ws.Range("E88").Value = 64

$ws.Range("B94").Value = 3302
$ws.Range("C94").Value = 15
$ws.Range("E94").Value = 1738

$ws.Range("B99").Value = 2593
$ws.Range("C99").Value = 135
$ws.Range("D99").Value = 1016
$ws.Range("E99").Value = 1569

$ws.Range("B107").Value = 1980
$ws.Range("C107").Value = 29
$ws.Range("D107").Value = 1548
$ws.Range("E107").Value = 421
$ws.Range("H107").Value = 11

$ws.Range("B108").Value = 1978
$ws.Range("C108").Value = 17
$ws.Range("D108").Value = 1302
$ws.Range("E108").Value = 565
$ws.Range("H108").Value = 111

$ws.Range("B128").Value = 1047
$ws.Range("C128").Value = 5
$ws.Range("D128").Value = 772
$ws.Range("E128").Value = 266
$ws.Range("H128").Value = 9

$ws.Range("B129").Value = 1046
$ws.Range("D129").Value = 913
$ws.Range("E129").Value = 66
$ws.Range("H129").Value = 67

$ws.Range("B130").Value = 990
$ws.Range("C130").Value = 2
$ws.Range("E130").Value = 147

$ws.Range("B142").Value = 757
$ws.Range("C142").Value = 20
$ws.Range("D142").Value = 206
$ws.Range("E142").Value = 546

$ws.Range("B156").Value = 378
$ws.Range("C156").Value = 11
$ws.Range("E156").Value = 54

$ws.Range("E160").Value = 178
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 9
